$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:G1) text -----------------------------------
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- Shift/replace the D,E,F numeric columns for rows 2-6 -------------
# New D = old E, new E = old F, new F = freshly supplied value.
$ws.Range("D2").Value = 14.22168333333333
$ws.Range("E2").Value = 209.93916
$ws.Range("F2").Value = 0.00039653827

$ws.Range("D3").Value = 11.27025133333333
$ws.Range("E3").Value = 175.49941
$ws.Range("F3").Value = 0.0003142445

$ws.Range("D4").Value = 8.521267333333334
$ws.Range("E4").Value = 110.15952
$ws.Range("F4").Value = 0.00023759554

$ws.Range("D5").Value = 30.16162866666667
$ws.Range("E5").Value = 390.83567
$ws.Range("F5").Value = 0.00084098624

$ws.Range("D6").Value = 12.74820333333333
$ws.Range("E6").Value = 229.36719
$ws.Range("F6").Value = 0.00035545374

# --- Add header-row cell comments (also creates the legacy VML drawing) -
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
